$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 3065.5715
$ws.Range("J17").Value = 3409.8333
$ws.Range("L17").Value = 10229.4999
$ws.Range("N17").Value = -10565.4999
# row 113
$ws.Range("H113").Value = 8200
$ws.Range("I113").Value = 8200
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 8200
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -4946
# row 131
$ws.Range("H131").Value = 3550
$ws.Range("J131").Value = 4333.3335
$ws.Range("L131").Value = 13000.0005
$ws.Range("N131").Value = -23080.0005
# row 135
$ws.Range("H135").Value = 1885.1428
$ws.Range("I135").Value = 639.4
$ws.Range("J135").Value = 4999.5
$ws.Range("K135").Value = 5754.599999999999
$ws.Range("L135").Value = 44995.5
$ws.Range("M135").Value = -3219.599999999999
$ws.Range("N135").Value = -50065.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 6669.6665
$ws.Range("I2").Value = 6669.6665
$ws.Range("K2").Value = 6669.6665
$ws.Range("M2").Value = -6556.6665
# row 61
$ws.Range("H61").Value = 1900
$ws.Range("I61").Value = 1900
$ws.Range("K61").Value = 1900
$ws.Range("M61").Value = -1688
# row 74
$ws.Range("H74").Value = 21038.576
$ws.Range("I74").Value = 20500.166
$ws.Range("J74").Value = 27499.5
$ws.Range("K74").Value = 20500.166
$ws.Range("L74").Value = 27499.5
$ws.Range("M74").Value = -19626.166
$ws.Range("N74").Value = -29247.5
# row 77
$ws.Range("H77").Value = 21038.576
$ws.Range("I77").Value = 20500.166
$ws.Range("J77").Value = 27499.5
$ws.Range("K77").Value = 102500.83
$ws.Range("L77").Value = 137497.5
$ws.Range("M77").Value = -98132.83
$ws.Range("N77").Value = -146233.5
# row 116
$ws.Range("H116").Value = 6669.6665
$ws.Range("I116").Value = 6669.6665
$ws.Range("K116").Value = 6669.6665
$ws.Range("M116").Value = -4375.6665
# row 136
$ws.Range("H136").Value = 1900
$ws.Range("I136").Value = 1900
$ws.Range("K136").Value = 5700
$ws.Range("M136").Value = -3150

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 6669.6665
$ws.Range("I3").Value = 6669.6665
$ws.Range("K3").Value = 6669.6665
$ws.Range("M3").Value = -6555.6665
# row 94
$ws.Range("H94").Value = 2801.2856
$ws.Range("I94").Value = 3703
$ws.Range("K94").Value = 3703
$ws.Range("M94").Value = -3252
# row 134
$ws.Range("H134").Value = 3997.7144
$ws.Range("I134").Value = 3997.3333
$ws.Range("K134").Value = 11991.9999
$ws.Range("M134").Value = -9456.999899999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2413.2856
$ws.Range("I31").Value = 2365
$ws.Range("K31").Value = 2365
$ws.Range("M31").Value = -2070
# row 34
$ws.Range("H34").Value = 2413.2856
$ws.Range("I34").Value = 2365
$ws.Range("K34").Value = 2365
$ws.Range("M34").Value = -2163
# row 58
$ws.Range("H58").Value = 4874.8887
$ws.Range("I58").Value = 4874.8887
$ws.Range("K58").Value = 4874.8887
$ws.Range("M58").Value = -4671.8887
# row 132
$ws.Range("H132").Value = 1663.2778
$ws.Range("I132").Value = 1141.8334
$ws.Range("K132").Value = 3425.5002
$ws.Range("M132").Value = -895.5001999999999
# row 134
$ws.Range("H134").Value = 4433.0586
$ws.Range("I134").Value = 3890.5386
$ws.Range("J134").Value = 6196.25
$ws.Range("K134").Value = 11671.6158
$ws.Range("L134").Value = 18588.75
$ws.Range("M134").Value = -9136.6158
$ws.Range("N134").Value = -23658.75
# row 136
$ws.Range("H136").Value = 4874.8887
$ws.Range("I136").Value = 4874.8887
$ws.Range("K136").Value = 14624.6661
$ws.Range("M136").Value = -12074.6661

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 81
$ws.Range("H81").Value = 5999.857
$ws.Range("J81").Value = 5833.3335
$ws.Range("L81").Value = 17500.0005
$ws.Range("N81").Value = -19746.0005
# row 84
$ws.Range("H84").Value = 5999.857
$ws.Range("J84").Value = 5833.3335
$ws.Range("L84").Value = 52500.0015
$ws.Range("N84").Value = -63732.0015

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 262.45456
$ws.Range("I2").Value = 362.42856
$ws.Range("K2").Value = 362.42856
$ws.Range("M2").Value = -249.42856
# row 70
$ws.Range("H70").Value = 4000
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 4000
$ws.Range("N70").Value = -4540
# row 73
$ws.Range("H73").Value = 4000
$ws.Range("J73").Value = 4000
$ws.Range("N73").Value = -5872

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 32
$ws.Range("H32").Value = 2800
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
# row 46
$ws.Range("H46").Value = 2178
$ws.Range("J46").Value = 3840.25
$ws.Range("L46").Value = 3840.25
$ws.Range("N46").Value = -4216.25
# row 55
$ws.Range("H55").Value = 1460.9412
$ws.Range("I55").Value = 2348.4443
$ws.Range("K55").Value = 2348.4443
$ws.Range("M55").Value = -2175.4443
# row 61
$ws.Range("H61").Value = 8589.799999999999
$ws.Range("I61").Value = 5737.25
$ws.Range("K61").Value = 5737.25
$ws.Range("M61").Value = -5535.25
# row 82
$ws.Range("H82").Value = 2558.3333
$ws.Range("J82").Value = 2192.5
$ws.Range("L82").Value = 2192.5
$ws.Range("N82").Value = -2914.5
# row 85
$ws.Range("H85").Value = 2558.3333
$ws.Range("J85").Value = 2192.5
$ws.Range("L85").Value = 2192.5
$ws.Range("N85").Value = -4688.5
# row 113
$ws.Range("H113").Value = 8589.799999999999
$ws.Range("I113").Value = 5737.25
$ws.Range("K113").Value = 5737.25
$ws.Range("M113").Value = -3567.25
# row 132
$ws.Range("H132").Value = 2854.44
$ws.Range("I132").Value = 2458
$ws.Range("K132").Value = 7374
$ws.Range("M132").Value = -4844

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 7364.2856
$ws.Range("J62").Value = 9150
$ws.Range("L62").Value = 9150
$ws.Range("N62").Value = -10398
# row 65
$ws.Range("H65").Value = 7364.2856
$ws.Range("J65").Value = 9150
$ws.Range("L65").Value = 45750
$ws.Range("N65").Value = -51990
# row 100
$ws.Range("H100").Value = 1378.6842
$ws.Range("J100").Value = 2199.8333
$ws.Range("L100").Value = 4399.6666
$ws.Range("N100").Value = -5481.6666
# row 136
$ws.Range("H136").Value = 1537.5
$ws.Range("I136").Value = 1537.5
$ws.Range("K136").Value = 4612.5
$ws.Range("M136").Value = -2062.5
